$d = $word.ActiveDocument

# Locate the paragraph "-statutes could be truths" (it currently carries the
# hidden _GoBack bookmark at its very end, right before the paragraph mark).
$targetPara = $null
foreach ($para in $d.Paragraphs) {
    $t = $para.Range.Text
    $t = $t.TrimEnd([char]13, [char]7)
    if ($t -eq "-statutes could be truths") {
        $targetPara = $para
        break
    }
}
if ($targetPara -eq $null) {
    throw "Could not find paragraph '-statutes could be truths'"
}

# Position right after the existing text, before the paragraph mark.
$insertPos = $targetPara.Range.End - 1

$hasGoBack = $d.Bookmarks.Exists("_GoBack")

if ($hasGoBack) {
    # The _GoBack bookmark sits collapsed at that same position. Word keeps a
    # collapsed bookmark "stuck" to the text immediately to its left as more
    # text is appended through the bookmark's own Range, so we grow it
    # forward with the new final line of text first (no paragraph mark yet -
    # this keeps the bookmark a simple collapsed point instead of wrapping a
    # whole paragraph, which is what happens if a bookmark is rebuilt
    # collapsed at a position that sits immediately before a paragraph
    # mark).
    $bm = $d.Bookmarks.Item("_GoBack")
    $bmRange = $bm.Range
    $bmRange.InsertAfter("-arrogance could be pride and vv")
} else {
    $tail = $d.Range($insertPos, $insertPos)
    $tail.InsertAfter("-arrogance could be pride and vv")
}

# Now splice in the new middle paragraph ("-saints could be holy ones")
# plus the two paragraph breaks that separate the three lines, inserting
# right after the original "-statutes could be truths" text. Because this
# insertion happens before the (now further along) bookmark position, the
# bookmark shifts along with it and ends up collapsed at the end of the
# newly typed "-arrogance could be pride and vv" paragraph, exactly as in
# the source edit.
$mid = $d.Range($insertPos, $insertPos)
$mid.InsertAfter([char]13 + "-saints could be holy ones" + [char]13)
